# Apply the daily cryptos-list data refresh (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a D-column price cell while preserving it as
# literal text (the sheet stores prices as text, not numbers), exactly as
# Excel keeps a number-format-forced cell as text.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "42.213.01"
$ws.Range("E2").Value = "  -1.24%  "

$ws.Range("D3").Value = "2.241.05"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("E4").Value = "  +0.17%  "

Set-TextValue "D5" "246.36"
$ws.Range("E5").Value = "  -0.91%  "

Set-TextValue "D6" "0.628"
$ws.Range("E6").Value = "  -2.23%  "

Set-TextValue "D7" "74.30"
$ws.Range("E7").Value = "  -3.45%  "

$ws.Range("E8").Value = "  +0.15%  "

Set-TextValue "D9" "0.616"
$ws.Range("E9").Value = "  -5.30%  "

Set-TextValue "D10" "42.07"
$ws.Range("E10").Value = "  +5.86%  "

$ws.Range("E11").Value = "  -3.01%  "

Set-TextValue "D12" "7.15"
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("E13").Value = "  -1.47%  "

Set-TextValue "D14" "14.49"
$ws.Range("E14").Value = "  -3.31%  "

Set-TextValue "D15" "0.848"
$ws.Range("E15").Value = "  -1.72%  "

$ws.Range("D16").Value = "2.232.41"
$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").Value = "42.073.82"
$ws.Range("E17").Value = "  -1.32%  "

$ws.Range("D18").Value = "0.0₃0983"
$ws.Range("E18").Value = "  -0.69%  "

Set-TextValue "D19" "72.15"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("E20").Value = "  -0.98%  "

Set-TextValue "D21" "231.32"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("E22").Value = "  +3.68%  "

Set-TextValue "D23" "8.71"
$ws.Range("E23").Value = "  +36.77%  "

$ws.Range("E24").Value = "  +0.07%  "

Set-TextValue "D25" "11.47"
$ws.Range("E25").Value = "  +0.92%  "

$ws.Range("E26").Value = "  -4.82%  "

$ws.Range("E27").Value = "  -1.78%  "

$ws.Range("E28").Value = "  +1.84%  "

Set-TextValue "D29" "169.20"
$ws.Range("E29").Value = "  +0.83%  "

Set-TextValue "D30" "20.63"
$ws.Range("E30").Value = "  -1.02%  "

Set-TextValue "D31" "0.0819"
$ws.Range("E31").Value = "  -3.57%  "

Set-TextValue "D32" "31.31"
$ws.Range("E32").Value = "  +3.38%  "

$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("E35").Value = "  +10.50%  "

$ws.Range("E36").Value = "  -1.50%  "

Set-TextValue "D37" "0.0313"
$ws.Range("E37").Value = "  +3.19%  "

Set-TextValue "D38" "13.69"

$ws.Range("E39").Value = "  -3.03%  "

$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D41" "0.205"
$ws.Range("E41").Value = "  -0.90%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D42" "62.13"
$ws.Range("E42").Value = "  +2.08%  "

Set-TextValue "D43" "106.34"
$ws.Range("E43").Value = "  -3.72%  "

Set-TextValue "D44" "0.102"
$ws.Range("E44").Value = "  +2.01%  "

Set-TextValue "D45" "8.65"
$ws.Range("E45").Value = "  -2.04%  "

$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("E47").Value = "  -2.73%  "

$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "2.28"
$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D50" "4.18"
$ws.Range("E50").Value = "  -6.30%  "

$ws.Range("E51").Value = "  +0.49%  "
